# Apply crypto price/volume refresh (GitHub Actions run, Fri Jun 14 16:54:28 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.421.84"
$ws.Range("E2").Value = "  -2.03%  "

$ws.Range("D3").Value = "3.407.99"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'593.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'142.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.409.01"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("D9").Value = "'0.467"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.26%  "

$ws.Range("D10").Value = "'0.135"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.98%  "

$ws.Range("D11").Value = "'7.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.96%  "

$ws.Range("D12").Value = "'0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.98%  "

$ws.Range("D13").Value = "3.983.31"
$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("E14").Value = "  -5.95%  "

$ws.Range("D15").Value = "'29.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.88%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.421.05"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "65.923.60"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").Value = "'10.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "

$ws.Range("D20").Value = "'6.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.70%  "

$ws.Range("D21").Value = "'14.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.90%  "

$ws.Range("D22").Value = "'417.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.85%  "

$ws.Range("D23").Value = "'0.582"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.44%  "

$ws.Range("D24").Value = "'77.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").Value = "3.541.91"
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").Value = "'0.0000110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.19%  "

$ws.Range("D28").Value = "'9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.56%  "

$ws.Range("D29").Value = "'7.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.82%  "

$ws.Range("D30").Value = "'2.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.25%  "

$ws.Range("E32").Value = "  -3.11%  "

$ws.Range("E33").Value = "  -8.46%  "

$ws.Range("D34").Value = "'24.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'1.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.71%  "

$ws.Range("D37").Value = "'5.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.22%  "

$ws.Range("D38").Value = "'7.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'172.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").Value = "'0.0865"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.35%  "

$ws.Range("D42").Value = "'5.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.81%  "

$ws.Range("D43").Value = "'0.871"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("E44").Value = "  -11.79%  "

$ws.Range("D45").Value = "'45.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").Value = "'26.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.08%  "

$ws.Range("D47").Value = "'1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.41%  "

$ws.Range("D48").Value = "'7.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.06%  "

$ws.Range("D49").Value = "'2.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.83%  "

$ws.Range("D50").Value = "'0.920"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.72%  "

$ws.Range("D51").Value = "'0.233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.78%  "
